# Apply the "Add dcr sheet" edit to the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "dcr" worksheet at the end of the workbook.
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "dcr"

# ---------------------------------------------------------------------
# 2. Populate the "dcr" sheet with header + year/value table.
# ---------------------------------------------------------------------
$newSheet.Range("A1").Value = "Stf"
$newSheet.Range("B1").Value = "Value"

$years = 2024..2050
$values = @(
    0.02,
    0.0222,
    0.0244,
    0.0266,
    0.0288,
    0.031,
    0.0332,
    0.0354,
    0.0376,
    0.0398,
    0.0398,
    0.042,
    0.0442,
    0.0464,
    0.0486,
    0.0508,
    0.053,
    0.0552,
    0.0574,
    0.0596,
    0.0618,
    0.064,
    0.0662,
    0.0684,
    0.0706,
    0.0728,
    0.075
)

for ($i = 0; $i -lt $years.Count; $i++) {
    $row = $i + 2
    $newSheet.Cells.Item($row, 1).Value = $years[$i]
    $newSheet.Cells.Item($row, 2).Value = $values[$i]
}

# Rows 3-27 (years 2025-2049) get a vertical-centered, wrap-text style.
# Build the combined style once on a scratch cell, then paste the
# resulting format onto the target range in a single operation so we
# don't leave an orphaned intermediate cellXf behind (each direct
# property write on a multi-cell range mints its own style record).
$stage = $newSheet.Range("Z1")
$stage.VerticalAlignment = -4108
$stage.WrapText = $true

$fmtRange = $newSheet.Range("B3:B27")
$stage.Copy()
$fmtRange.PasteSpecial(-4122)
$stage.Clear()
$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 3. Update the Params sheet: "anti duping Index" (B13) 1 -> 0
# ---------------------------------------------------------------------
$paramsSheet = $wb.Worksheets.Item("Params")
$paramsSheet.Range("B13").Value = 0
$paramsSheet.Range("B18").Select()

# ---------------------------------------------------------------------
# 4. Update the eu_secondary_cost sheet selection (whole-column select,
#    no longer the active/tabSelected sheet).
# ---------------------------------------------------------------------
$euSecondary = $wb.Worksheets.Item("eu_secondary_cost")
$euSecondary.Columns("A:B").Select()

# ---------------------------------------------------------------------
# 5. Activate the new "dcr" sheet (becomes tabSelected) with its
#    own selection.
# ---------------------------------------------------------------------
$newSheet.Activate()
$newSheet.Range("I19").Select()
